$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move column F (which held the "Cargo" header/data) one column to the left
# into E, since column E is currently empty. Copy cell-by-cell (via Formula,
# which round-trips plain text faithfully) to preserve per-row values, then
# clear the now-vacated F column.
for ($r = 1; $r -le 23; $r++) {
    $srcCell = $ws.Cells.Item($r, 6)
    $dstCell = $ws.Cells.Item($r, 5)
    $dstCell.Formula = $srcCell.Formula
    $srcCell.Clear()
}

# The header cell (row 1) was bold in column F; re-apply that formatting on
# its new home in column E so the moved header keeps the same look.
$ws.Cells.Item(1, 5).Font.Bold = $true

# Update the selection to mirror the new active selection on column E
# (E1 active cell, full-column selection E1:E1048576).
$ws.Range("E1:E1048576").Select()
